# Add models Item, RedressKit and RedressKitConsist
# (new "redress" reference table placed in column F next to "Service Level",
#  and a new "steps" reference table placed in column D next to "Category")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- remove stray empty placeholder cell in the "tool" block ---
$ws.Range("D10").ClearContents()

# --- new "redress" model header, column F, next to "Service Level" ---
$ws.Range("F26").Value = "redress"
$ws.Range("F26").Font.Bold = $true

# --- "redress" model fields, column F ---
$ws.Range("F27").Value = "equipment_id"
$ws.Range("F27").Borders.LineStyle = 1

$ws.Range("F28").Value = "step_id"
$ws.Range("F28").Borders.LineStyle = 1

$ws.Range("F29").Value = "position"
$ws.Range("F29").Borders.LineStyle = 1

# --- new "steps" model header, column D, reusing the bold header look ---
$ws.Range("D28").Value = "steps"
$ws.Range("D28").Font.Bold = $true

# --- "steps" model fields, column D ---
$ws.Range("D29").Value = "id"
$ws.Range("D29").Borders.LineStyle = 1

$ws.Range("D30").Value = "description"
$ws.Range("D30").Borders.LineStyle = 1

$ws.Range("D31").Value = "level_id"
$ws.Range("D31").Borders.LineStyle = 1

# --- drop the now-obsolete empty placeholder cells under "Category" ---
$ws.Range("D32").ClearContents()
$ws.Range("D33").ClearContents()

# --- restore the selection left by the author at the end of the edit ---
$ws.Range("F26").Select()
